# Fruta / hortaliza, semanal
#
# A new weekly price-report group (date 2023-12-21 / serial 45281) is
# inserted right above the existing 2023-11-06 (serial 45236) group at
# rows 60-63. That new group carries the same Especial/Primera/Segunda/
# Tercera values that the 2023-11-06 group already had, so the simplest
# faithful reproduction is: insert 4 blank rows above the current
# 60-63 block (which pushes that block - and everything below it - down
# by 4 rows), then fill the vacated 60-63 rows with a copy of what is
# now sitting in rows 64-67 (the original 2023-11-06 data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("60:63").Insert()

$ws.Range("A64:T67").Copy()
$ws.Range("A60").PasteSpecial()
